$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values for rows 2-5 (new 1000-point dataset slice) ---
$arr = New-Object 'object[,]' 4,34
$arr[0,0] = 45102.50694444445
$arr[0,1] = 20.658
$arr[0,2] = 14.028
$arr[0,3] = 4.08
$arr[0,4] = 43.704
$arr[0,5] = 35.637
$arr[0,6] = 16.257
$arr[0,7] = 52.557
$arr[0,8] = 25.014
$arr[0,9] = 10.522
$arr[0,10] = 16.084
$arr[0,11] = 17.263
$arr[0,12] = 18.024
$arr[0,13] = 5.19
$arr[0,14] = 16.166
$arr[0,15] = 22.644
$arr[0,16] = 13.689
$arr[0,17] = 3.709
$arr[0,18] = 2.453
$arr[0,19] = 238.72
$arr[0,20] = 44.979
$arr[0,21] = 14.922
$arr[0,22] = 29.658
$arr[0,23] = 15.316
$arr[0,24] = 2.936
$arr[0,25] = 26.039
$arr[0,26] = 13.181
$arr[0,27] = 11.919
$arr[0,28] = 13.951
$arr[0,29] = 17.711
$arr[0,30] = 3.458
$arr[0,31] = 46.546
$arr[0,32] = 8.338
$arr[0,33] = 18.656
$arr[1,0] = 45102.51388888889
$arr[1,1] = 23.541
$arr[1,2] = 16.982
$arr[1,3] = 2.047
$arr[1,4] = 50.805
$arr[1,5] = 41.68
$arr[1,6] = 18.525
$arr[1,7] = 70.878
$arr[1,8] = 28.505
$arr[1,9] = 12.477
$arr[1,10] = 18.64
$arr[1,11] = 20.381
$arr[1,12] = 21.369
$arr[1,13] = 5.918
$arr[1,14] = 18.422
$arr[1,15] = 26.097
$arr[1,16] = 15.636
$arr[1,17] = 1.663
$arr[1,18] = 1.291
$arr[1,19] = 273.097
$arr[1,20] = 51.526
$arr[1,21] = 17.004
$arr[1,22] = 34.404
$arr[1,23] = 18.101
$arr[1,24] = 2.921
$arr[1,25] = 34.387
$arr[1,26] = 15.02
$arr[1,27] = 13.45
$arr[1,28] = 15.781
$arr[1,29] = 21.162
$arr[1,30] = 1.266
$arr[1,31] = 64.353
$arr[1,32] = 9.562
$arr[1,33] = 21.259
$arr[2,0] = 45102.52083333334
$arr[2,1] = 6.726
$arr[2,2] = 4.593
$arr[2,3] = 0.996
$arr[2,4] = 14.404
$arr[2,5] = 11.642
$arr[2,6] = 5.294
$arr[2,7] = 26.097
$arr[2,8] = 8.144
$arr[2,9] = 3.481
$arr[2,10] = 5.085
$arr[2,11] = 5.813
$arr[2,12] = 5.997
$arr[2,13] = 1.698
$arr[2,14] = 5.263
$arr[2,15] = 7.411
$arr[2,16] = 4.668
$arr[2,17] = 0.991
$arr[2,18] = 0.536
$arr[2,19] = 72.832
$arr[2,20] = 14.975
$arr[2,21] = 4.858
$arr[2,22] = 9.811
$arr[2,23] = 5.142
$arr[2,24] = 1.056
$arr[2,25] = 11.815
$arr[2,26] = 4.291
$arr[2,27] = 3.964
$arr[2,28] = 4.629
$arr[2,29] = 5.936
$arr[2,30] = 0.773
$arr[2,31] = 23.997
$arr[2,32] = 2.658
$arr[2,33] = 6.076
$arr[3,0] = 45102.52777777778
$arr[3,1] = 13.45
$arr[3,2] = 9.76
$arr[3,3] = 1.01
$arr[3,4] = 29.11
$arr[3,5] = 23.84
$arr[3,6] = 10.59
$arr[3,7] = 39.71
$arr[3,8] = 16.29
$arr[3,9] = 7.15
$arr[3,10] = 10.62
$arr[3,11] = 11.72
$arr[3,12] = 12.28
$arr[3,13] = 3.38
$arr[3,14] = 10.53
$arr[3,15] = 14.9
$arr[3,16] = 8.97
$arr[3,17] = 0.8
$arr[3,18] = 0.63
$arr[3,19] = 152.89
$arr[3,20] = 29.37
$arr[3,21] = 9.72
$arr[3,22] = 19.61
$arr[3,23] = 10.37
$arr[3,24] = 1.67
$arr[3,25] = 19.17
$arr[3,26] = 8.58
$arr[3,27] = 7.68
$arr[3,28] = 9.01
$arr[3,29] = 12.2
$arr[3,30] = 0.56
$arr[3,31] = 35.77
$arr[3,32] = 5.45
$arr[3,33] = 12.15
$ws.Range("A2:AH5").Value = $arr

# --- Remove the now-unused trailing row 6 ---
$ws.Rows(6).Delete()

# --- Widen most data columns from 7 (and the narrower col L=6) to 8 chars; col T to 9 ---
$ws.Range("B:C").ColumnWidth = 7.17
$ws.Range("E:M").ColumnWidth = 7.17
$ws.Range("O:Q").ColumnWidth = 7.17
$ws.Range("U:X").ColumnWidth = 7.17
$ws.Range("Z:AD").ColumnWidth = 7.17
$ws.Range("AF:AF").ColumnWidth = 7.17
$ws.Range("AH:AH").ColumnWidth = 7.17
$ws.Range("T:T").ColumnWidth = 8.17
